$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellAddr, $NewValue)
    $r = $ws.Range($CellAddr)
    # Leading apostrophe forces Excel to store the value as literal text,
    # matching the inlineStr cell type used throughout this sheet, even
    # though the text looks like a number.
    $r.Formula = "'" + $NewValue
    # Re-apply the default style so the quote-prefix formatting introduced
    # by the apostrophe entry does not leave a visible style change behind.
    $r.Style = "Normal"
}

Set-TextValue "D2" "274.80"
Set-TextValue "D3" "22.93"
Set-TextValue "D4" "6.426"
Set-TextValue "D5" "0.06278"
Set-TextValue "D6" "3.661"
Set-TextValue "D7" "6.670"
Set-TextValue "D8" "1.401"
Set-TextValue "D9" "0.8307"
Set-TextValue "D10" "0.01380"
Set-TextValue "D11" "0.1625"
Set-TextValue "D12" "0.08304"
Set-TextValue "D14" "0.03096"
Set-TextValue "D15" "0.09298"
Set-TextValue "D16" "3.880"
Set-TextValue "D17" "0.001655"
Set-TextValue "D18" "0.04782"
Set-TextValue "D19" "0.006415"
Set-TextValue "D21" "0.001091"
Set-TextValue "D24" "2.387"
Set-TextValue "D25" "0.3348"
Set-TextValue "D40" "0.04712"
Set-TextValue "D41" "0.007049"
Set-TextValue "D43" "0.003700"
Set-TextValue "D44" "0.01185"
Set-TextValue "D45" "0.00006249"
Set-TextValue "D48" "0.7965"
Set-TextValue "D49" "0.04008"
Set-TextValue "D50" "0.00002300"
